$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn phụ phẫu 2")

# New totals row (row 2) appended below the header row.
# Text columns (A, C, D, E, F, G, H, J, Q, R, S, T) stay blank.
# "Mã dịch vụ" (B) is a numeric column left blank (no value).
# The remaining numeric/money columns (I, K, L, M, N, O, P) are zeroed.
$ws.Cells.Item(2, 1).Value = ""
$ws.Cells.Item(2, 2).Value = $null
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = ""
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = ""
$ws.Cells.Item(2, 18).Value = ""
$ws.Cells.Item(2, 19).Value = ""
$ws.Cells.Item(2, 20).Value = ""
